$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.150.43"
$ws.Range("E2").Value = "  +4.91%  "
$ws.Range("D3").Value = "3.361.83"
$ws.Range("E3").Value = "  +10.36%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'257.29"
$ws.Range("E5").Value = "  +10.26%  "
$ws.Range("D6").Value = "'622.23"
$ws.Range("E6").Value = "  +2.70%  "
$ws.Range("E7").Value = "  +10.61%  "
$ws.Range("E8").Value = "  +2.49%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "3.358.53"
$ws.Range("E10").Value = "  +10.44%  "
$ws.Range("D11").Value = "'0.804"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("E12").Value = "  +2.32%  "
$ws.Range("D13").Value = "97.918.63"
$ws.Range("E13").Value = "  +4.81%  "
$ws.Range("D14").Value = "'35.72"
$ws.Range("E14").Value = "  +6.91%  "
$ws.Range("D15").Value = "'0.0000245"
$ws.Range("E15").Value = "  +2.79%  "
$ws.Range("D16").Value = "3.962.76"
$ws.Range("E16").Value = "  +9.54%  "
$ws.Range("E17").Value = "  +4.60%  "
$ws.Range("D18").Value = "3.355.98"
$ws.Range("E18").Value = "  +9.77%  "
$ws.Range("D19").Value = "'3.60"
$ws.Range("E19").Value = "  +2.35%  "
$ws.Range("E20").Value = "  +4.91%  "
$ws.Range("D21").Value = "'483.32"
$ws.Range("E21").Value = "  +10.93%  "
$ws.Range("D22").Value = "'5.83"
$ws.Range("E22").Value = "  +3.17%  "
$ws.Range("D23").Value = "'0.0000204"
$ws.Range("E23").Value = "  +8.23%  "
$ws.Range("D24").Value = "'9.13"
$ws.Range("E24").Value = "  +4.49%  "
$ws.Range("E25").Value = "  +3.53%  "
$ws.Range("D26").Value = "'88.29"
$ws.Range("E26").Value = "  +4.64%  "
$ws.Range("D27").Value = "'12.01"
$ws.Range("E27").Value = "  +2.93%  "
$ws.Range("E28").Value = "  +10.33%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "'0.250"
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("E31").Value = "  +4.78%  "
$ws.Range("E32").Value = "  -8.62%  "
$ws.Range("D33").Value = "'0.122"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").Value = "'9.24"
$ws.Range("E34").Value = "  +2.72%  "
$ws.Range("D35").Value = "'27.22"
$ws.Range("E35").Value = "  +7.95%  "
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "'517.07"
$ws.Range("E37").Value = "  +12.68%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.151"
$ws.Range("E38").Value = "  -2.30%  "
$ws.Range("E39").Value = "  +2.68%  "
$ws.Range("E40").Value = "  +4.00%  "
$ws.Range("D41").Value = "'0.446"
$ws.Range("E41").Value = "  +2.30%  "
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("E43").Value = "  -4.69%  "
$ws.Range("D44").Value = "'3.24"
$ws.Range("E44").Value = "  +4.48%  "
$ws.Range("D46").Value = "'0.777"
$ws.Range("E46").Value = "  +17.32%  "
$ws.Range("D47").Value = "'160.78"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("E48").Value = "  +5.37%  "
$ws.Range("E49").Value = "  +7.92%  "
$ws.Range("D50").Value = "'45.54"
$ws.Range("E50").Value = "  +4.17%  "
$ws.Range("E51").Value = "  +7.44%  "
